$wb = $excel.ActiveWorkbook

# The localization status report moved from "Ready for handoff" to
# "In Translation" for every language status cell (Overview!E:F and the
# per-language "Status" column on each language sheet).
foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    foreach ($cell in $usedRange.Cells) {
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value2 = "In Translation"
        }
    }
}

# "In Translation" is shorter than "Ready for handoff", so the status
# columns shrink to fit the new text (Excel's column best-fit behavior).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
